$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "'51.840.05"
$ws.Cells.Item(2,5).Value = "'  +0.15%  "
$ws.Cells.Item(3,4).Value = "'2.837.32"
$ws.Cells.Item(3,5).Value = "'  +2.73%  "
$ws.Cells.Item(4,5).Value = "'  -0.04%  "
$ws.Cells.Item(5,4).Value = "'352.55"
$ws.Cells.Item(5,5).Value = "'  +5.82%  "
$ws.Cells.Item(6,4).Value = "'113.97"
$ws.Cells.Item(6,5).Value = "'  -2.13%  "
$ws.Cells.Item(7,4).Value = "'0.560"
$ws.Cells.Item(7,5).Value = "'  +3.77%  "
$ws.Cells.Item(8,5).Value = "'  -0.08%  "
$ws.Cells.Item(9,4).Value = "'0.601"
$ws.Cells.Item(9,5).Value = "'  +4.14%  "
$ws.Cells.Item(10,4).Value = "'41.71"
$ws.Cells.Item(10,5).Value = "'  -0.27%  "
$ws.Cells.Item(11,4).Value = "'0.0853"
$ws.Cells.Item(11,5).Value = "'  -0.77%  "
$ws.Cells.Item(12,4).Value = "'20.00"
$ws.Cells.Item(12,5).Value = "'  -1.01%  "
$ws.Cells.Item(13,5).Value = "'  +1.53%  "
$ws.Cells.Item(14,4).Value = "'7.74"
$ws.Cells.Item(14,5).Value = "'  +1.09%  "
$ws.Cells.Item(15,4).Value = "'3.275.86"
$ws.Cells.Item(15,5).Value = "'  +2.54%  "
$ws.Cells.Item(16,4).Value = "'2.828.28"
$ws.Cells.Item(16,5).Value = "'  +1.76%  "
$ws.Cells.Item(17,4).Value = "'0.898"
$ws.Cells.Item(17,5).Value = "'  +0.92%  "
$ws.Cells.Item(18,4).Value = "'51.692.41"
$ws.Cells.Item(18,5).Value = "'  -0.02%  "
$ws.Cells.Item(19,5).Value = "'  +7.62%  "
$ws.Cells.Item(20,5).Value = "'  -1.97%  "
$ws.Cells.Item(21,4).Value = "'13.49"
$ws.Cells.Item(21,5).Value = "'  -0.16%  "
$ws.Cells.Item(22,4).Value = "'0.0₃0995"
$ws.Cells.Item(22,5).Value = "'  +2.15%  "
$ws.Cells.Item(23,4).Value = "'270.99"
$ws.Cells.Item(23,5).Value = "'  -2.69%  "
$ws.Cells.Item(24,4).Value = "'69.76"
$ws.Cells.Item(24,5).Value = "'  +0.18%  "
$ws.Cells.Item(25,5).Value = "'  +3.48%  "
$ws.Cells.Item(26,4).Value = "'26.75"
$ws.Cells.Item(26,5).Value = "'  -0.12%  "
$ws.Cells.Item(27,5).Value = "'  +0.00%  "
$ws.Cells.Item(28,4).Value = "'10.31"
$ws.Cells.Item(28,5).Value = "'  +1.28%  "
$ws.Cells.Item(29,5).Value = "'  +1.27%  "
$ws.Cells.Item(30,5).Value = "'  -1.40%  "
$ws.Cells.Item(31,4).Value = "'50.71"
$ws.Cells.Item(31,5).Value = "'  +1.36%  "
$ws.Cells.Item(32,4).Value = "'33.89"
$ws.Cells.Item(32,5).Value = "'  -3.34%  "
$ws.Cells.Item(33,4).Value = "'0.0450"
$ws.Cells.Item(33,5).Value = "'  +27.94%  "
$ws.Cells.Item(34,4).Value = "'5.82"
$ws.Cells.Item(34,5).Value = "'  +4.45%  "
$ws.Cells.Item(35,4).Value = "'0.0828"
$ws.Cells.Item(35,5).Value = "'  +0.51%  "
$ws.Cells.Item(36,4).Value = "'0.999"
$ws.Cells.Item(36,5).Value = "'  -0.14%  "
$ws.Cells.Item(37,5).Value = "'  -0.13%  "
$ws.Cells.Item(38,4).Value = "'3.23"
$ws.Cells.Item(38,5).Value = "'  -0.15%  "
$ws.Cells.Item(39,4).Value = "'4.89"
$ws.Cells.Item(39,5).Value = "'  -2.44%  "
$ws.Cells.Item(40,4).Value = "'18.03"
$ws.Cells.Item(40,5).Value = "'  -4.94%  "
$ws.Cells.Item(41,2).Value = "Stacks"
$ws.Cells.Item(41,3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(41,4).Value = "'2.56"
$ws.Cells.Item(41,5).Value = "'  +5.09%  "
$ws.Cells.Item(42,2).Value = "EnergySwap"
$ws.Cells.Item(42,3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(42,4).Value = "'23.52"
$ws.Cells.Item(42,5).Value = "'  +1.19%  "
$ws.Cells.Item(43,5).Value = "'  +1.14%  "
$ws.Cells.Item(44,4).Value = "'125.61"
$ws.Cells.Item(44,5).Value = "'  -1.12%  "
$ws.Cells.Item(45,5).Value = "'  +0.25%  "
$ws.Cells.Item(46,4).Value = "'2.083.01"
$ws.Cells.Item(46,5).Value = "'  -0.28%  "
$ws.Cells.Item(47,4).Value = "'3.35"
$ws.Cells.Item(47,5).Value = "'  +0.83%  "
$ws.Cells.Item(48,5).Value = "'  +3.58%  "
$ws.Cells.Item(49,4).Value = "'5.71"
$ws.Cells.Item(49,5).Value = "'  +3.01%  "
$ws.Cells.Item(50,4).Value = "'0.939"
$ws.Cells.Item(50,5).Value = "'  +7.22%  "
$ws.Cells.Item(51,4).Value = "'60.77"
$ws.Cells.Item(51,5).Value = "'  +1.38%  "
